$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant (used to copy a clean "General"-text cell style
# onto cells that change from a number to the "0"/"***.*" placeholder text,
# so the result matches the existing placeholder-cell styling exactly).
$xlPasteFormats = -4122

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/30/2023  Through  2/5/2023"

# --- Data table updates (rows 14-29) ---

# Row 14: M14 goes from the text placeholder "***.*" to numeric 0,
# matching the percent-change number format used by L14/N14.
$ws.Range("M14").NumberFormat = $ws.Range("L14").NumberFormat
$ws.Range("M14").Value = 0

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("D15").PasteSpecial($xlPasteFormats)
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E15").PasteSpecial($xlPasteFormats)
$ws.Range("M15").Value = -66.666666666666
$ws.Range("N15").Value = -92.307692307692

# Row 16
$ws.Range("C16").Value = 2
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("D16").PasteSpecial($xlPasteFormats)
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E16").PasteSpecial($xlPasteFormats)
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 125
$ws.Range("I16").Value = 26
$ws.Range("K16").Value = 160
$ws.Range("L16").Value = 160
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = -78.333333333333

# Row 17
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 31
$ws.Range("H17").Value = 19.230769230769
$ws.Range("I17").Value = 36
$ws.Range("J17").Value = 33
$ws.Range("K17").Value = 9.090909090909
$ws.Range("L17").Value = 89.473684210526
$ws.Range("M17").Value = 71.428571428571
$ws.Range("N17").Value = -47.826086956521

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("D18").PasteSpecial($xlPasteFormats)
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E18").PasteSpecial($xlPasteFormats)
$ws.Range("F18").Value = 13
$ws.Range("H18").Value = 225
$ws.Range("I18").Value = 15
$ws.Range("K18").Value = 200
$ws.Range("L18").Value = 66.666666666666
$ws.Range("M18").Value = -40
$ws.Range("N18").Value = -89.285714285714

# Row 19
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = 5.714285714285
$ws.Range("I19").Value = 46
$ws.Range("J19").Value = 44
$ws.Range("K19").Value = 4.545454545454
$ws.Range("L19").Value = 170.588235294118
$ws.Range("M19").Value = 48.387096774193
$ws.Range("N19").Value = -17.857142857142

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 300
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 175
$ws.Range("I20").Value = 12
$ws.Range("J20").Value = 4
$ws.Range("K20").Value = 200
$ws.Range("L20").Value = 140
$ws.Range("M20").Value = -40
$ws.Range("N20").Value = -92.727272727272

# Row 21
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 71.428571428571
$ws.Range("F21").Value = 112
$ws.Range("G21").Value = 80
$ws.Range("H21").Value = 40
$ws.Range("I21").Value = 137
$ws.Range("J21").Value = 99
$ws.Range("K21").Value = 38.383838383838
$ws.Range("L21").Value = 117.460317460317
$ws.Range("M21").Value = 7.874015748031
$ws.Range("N21").Value = -75.709219858156

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("D23").PasteSpecial($xlPasteFormats)
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E23").PasteSpecial($xlPasteFormats)
$ws.Range("F23").Value = 12
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 200
$ws.Range("I23").Value = 13
$ws.Range("K23").Value = 160
$ws.Range("L23").Value = 225
$ws.Range("M23").Value = 550

# Row 24
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 110
$ws.Range("F24").Value = 97
$ws.Range("G24").Value = 60
$ws.Range("H24").Value = 61.666666666666
$ws.Range("I24").Value = 123
$ws.Range("J24").Value = 75
$ws.Range("K24").Value = 64
$ws.Range("L24").Value = 41.379310344827
$ws.Range("M24").Value = -11.510791366906

# Row 25
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = -15.384615384615
$ws.Range("F25").Value = 55
$ws.Range("G25").Value = 46
$ws.Range("H25").Value = 19.565217391304
$ws.Range("I25").Value = 66
$ws.Range("J25").Value = 57
$ws.Range("K25").Value = 15.78947368421
$ws.Range("L25").Value = 106.25
$ws.Range("M25").Value = -38.317757009345

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("D26").PasteSpecial($xlPasteFormats)
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E26").PasteSpecial($xlPasteFormats)
$ws.Range("L26").Value = -66.666666666666

# Row 27
$ws.Range("C27").Value = 2
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("D27").PasteSpecial($xlPasteFormats)
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E27").PasteSpecial($xlPasteFormats)
$ws.Range("F27").Value = 6
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 7
$ws.Range("K27").Value = 75
$ws.Range("L27").Value = 75

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("D28").PasteSpecial($xlPasteFormats)
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E28").PasteSpecial($xlPasteFormats)

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("D29").PasteSpecial($xlPasteFormats)
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E29").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false
